$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bullet = [char]0x2022

# The "Patients Database" block (rows 105-113) is being condensed: two detail
# rows are dropped (9 rows -> 7 rows). Deleting the old "* Tab Tracking" (108)
# and "* Registration" (109) rows lets every row below shift up by two and
# land on a sibling that already carries the styling/formula shape we want
# to end up with (the old "* Meetings" row becomes the new "Other" row, the
# old totals row keeps its SUM formula shape, and the old @Parsiss/@Home
# rows keep their number formatting) - Excel also auto-adjusts the SUM()
# range and the I4 grand-total formula for us.
$ws.Rows.Item(109).Delete()
$ws.Rows.Item(108).Delete()

# --- Detail rows (105-108) -------------------------------------------------
$ws.Range("B105").Value = "* Patients Database"
$ws.Range("C105").Value = 9
$ws.Range("E105").Value = $bullet + " Patients Database"

$ws.Range("B106").Value = "* GUI"
$ws.Range("C106").Value = 6

$ws.Range("B107").Value = "* Registration"
$ws.Range("C107").Value = 15

$ws.Range("B108").Value = "* Other"
$ws.Range("C108").Value = 5

# Row 109 (totals) already reads "Total Hours" / =SUM(C105:C108) - it is the
# old totals row (111) shifted up, with its label/SUM formula auto-adjusted
# by the row delete above - nothing to change there.

# --- Paid / Not paid rows (110-111) -----------------------------------------
# C110/C111 already read "@Parsiss"/"@Home" (inherited from the old rows
# 112/113 that shifted into these slots) - leave them untouched so they keep
# their original shared-string + quote-prefixed-text cell style.
$ws.Range("D110").Value = 1
$ws.Range("D111").Formula = "=C109-D110"

# Match the author's final cursor position/selection.
$ws.Range("D112").Select() | Out-Null
